$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'27.656.71"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Formula = "'1.871.91"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("D4").Formula = "'1.004"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Formula = "'332.11"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +2.59%  '
$ws.Range("D6").Formula = "'1.004"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("D7").Formula = "'0.4723"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +3.84%  '
$ws.Range("D8").Formula = "'0.3945"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +1.81%  '
$ws.Range("D9").Formula = "'48.01"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").Formula = "'0.08067"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +1.82%  '
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Formula = "'22.01"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +2.75%  '
$ws.Range("D13").Formula = "'1.890.39"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("D14").Formula = "'5.959"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +0.83%  '
$ws.Range("D15").Formula = "'7.148"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").Formula = "'0.00001047"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +1.41%  '
$ws.Range("D18").Formula = "'86.85"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").Formula = "'0.06641"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").Formula = "'17.24"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").Formula = "'27.657.41"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("D23").Formula = "'5.509"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").Formula = "'10.99"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("D25").Formula = "'2.313"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("D26").Formula = "'2.112.30"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("D27").Formula = "'158.82"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +3.82%  '
$ws.Range("D28").Formula = "'20.27"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +2.37%  '
$ws.Range("D29").Formula = "'2.099"
$ws.Range("D29").Style = 'Normal'
$ws.Range("D30").Formula = "'5.568"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").Formula = "'122.23"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +1.08%  '
$ws.Range("D32").Formula = "'0.9708"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +3.82%  '
$ws.Range("D33").Formula = "'0.09540"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("D34").Formula = "'1.453"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -2.71%  '
$ws.Range("D35").Formula = "'3.592"
$ws.Range("D35").Style = 'Normal'
$ws.Range("D36").Formula = "'5.345"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +1.32%  '
$ws.Range("E37").Value = '  +1.64%  '
$ws.Range("D38").Formula = "'0.02258"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").Formula = "'1.229"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +0.80%  '
$ws.Range("D40").Formula = "'8.177"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("D41").Formula = "'0.6036"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +2.16%  '
$ws.Range("D42").Formula = "'0.1901"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +0.60%  '
$ws.Range("D43").Formula = "'10.28"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +1.34%  '
$ws.Range("D44").Formula = "'1.272"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").Formula = "'0.5697"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +1.22%  '
$ws.Range("D46").Formula = "'12.25"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +2.27%  '
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("D48").Formula = "'3.379"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("E50").Value = '  +5.88%  '
$ws.Range("D51").Formula = "'0.00000000301"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +15.72%  '
